$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - numeric values
$ws.Range("A2").Value = 5000
$ws.Range("B2").Value = 5000

# Row 3 - numeric values
$ws.Range("A3").Value = 6000
$ws.Range("B3").Value = 6000.6

# Row 4 - text values (stored as text/strings, not numbers)
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "6000.6"
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "5000.0"
$ws.Range("B4").Style = "Normal"
